# Update the 100 arithmetic-expression cells (20 rows x 5 columns) in the
# first table to the new set of expressions, preserving all existing
# run/paragraph formatting (only the visible text changes per cell).

$d = $word.ActiveDocument

$newValues = @("31+18=","34+7=","50+19=","36+21=","12+19=","13+31=","50+4=","53-6=","90+1=","88-26=","64-52=","10+10=","1+74=","15+6=","88-78=","98-13=","85-66=","65-15=","75-60=","37+2=","79-75=","56+38=","33+43=","93-65=","1+37=","28+10=","5+58=","87-27=","22-9=","73+25=","91-42=","1+63=","46-24=","45-24=","94-37=","49+23=","62+19=","67-24=","31+12=","95-43=","51-31=","47+48=","96-29=","4+0=","75-43=","90-47=","88-6=","31-21=","39-16=","90-67=","35+62=","21+14=","63-50=","88-13=","16+7=","59+26=","3+25=","91-11=","47+1=","0+32=","21+76=","4+63=","38-27=","56+11=","92-72=","90-80=","82-37=","24+33=","0+21=","66-29=","44-0=","81+2=","19-10=","37+51=","69-15=","76-55=","33+55=","89-29=","56-40=","80+13=","8+70=","67-6=","20-13=","62+11=","62-25=","0+57=","42-20=","22+15=","80-17=","77-56=","12+35=","4+38=","70-36=","95-1=","53-31=","73+6=","72-53=","82-63=","16+26=","6+16=")

$table = $d.Tables.Item(1)
$cols = $table.Columns.Count

$idx = 0
for ($row = 1; $row -le $table.Rows.Count; $row++) {
    for ($col = 1; $col -le $cols; $col++) {
        $cell = $table.Cell($row, $col)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output ("Updated " + $idx + " cells")
